$wb = $excel.ActiveWorkbook

# --- "Prix Spot" sheet: insert a new date column before the old "DI" column ---
# The source used to jump straight from 03-nov (DH) to 01-oct. (old DI).
# A new day, 04-nov, needs to be inserted right after 03-nov, which pushes every
# column from the old DI (01-oct.) through EM (31-oct.) one column to the right
# (new range DJ:EN), growing the sheet's used range from EM25 to EN25.
$ws = $wb.Worksheets.Item("Prix Spot")
$ws.Range("DI1:DI25").Insert(-4161)  # xlShiftToRight

# Populate the freshly inserted column: header + "-" placeholders for every hour row.
$ws.Range("DI1").Value = "04-nov"
$ws.Range("DI2:DI25").Value = "-"

# --- "Gaz" sheet: correct the forecasted price for 2025-11-01 / 2025-11-02 ---
$gaz = $wb.Worksheets.Item("Gaz")
$gaz.Range("B140").Value = 29.3
$gaz.Range("B141").Value = 29.3
